$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.465.53'
$ws.Range("E2").Value = '  +3.24%  '
$ws.Range("D3").Value = '2.313.89'
$ws.Range("E3").Value = '  +1.95%  '
$ws.Range("E4").Value = '  -0.15%  '
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = '516.78'
$r.ClearFormats()
$ws.Range("E5").Value = '  +2.59%  '
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = '134.93'
$r.ClearFormats()
$ws.Range("E6").Value = '  +5.78%  '
$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = '0.996'
$r.ClearFormats()
$ws.Range("E7").Value = '  -0.25%  '
$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = '0.536'
$r.ClearFormats()
$ws.Range("E8").Value = '  +1.51%  '
$ws.Range("D9").Value = '2.332.88'
$ws.Range("E10").Value = '  +4.91%  '
$ws.Range("E11").Value = '  -0.97%  '
$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = '5.34'
$r.ClearFormats()
$ws.Range("E12").Value = '  +5.11%  '
$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = '0.340'
$r.ClearFormats()
$ws.Range("E13").Value = '  +0.52%  '
$ws.Range("E14").Value = '  +2.14%  '
$ws.Range("D15").Value = '2.727.62'
$ws.Range("E15").Value = '  +1.99%  '
$ws.Range("D16").Value = '56.586.99'
$ws.Range("E16").Value = '  +3.25%  '
$ws.Range("E17").Value = '  +2.89%  '
$ws.Range("D18").Value = '2.319.31'
$ws.Range("E18").Value = '  +1.10%  '
$ws.Range("E19").Value = '  +1.46%  '
$ws.Range("E20").Value = '  +1.36%  '
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = '322.99'
$r.ClearFormats()
$ws.Range("E21").Value = '  +3.96%  '
$ws.Range("E22").Value = '  +0.16%  '
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = '0.999'
$r.ClearFormats()
$ws.Range("E23").Value = '  +0.21%  '
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = '60.69'
$r.ClearFormats()
$ws.Range("E24").Value = '  +1.59%  '
$ws.Range("E25").Value = '  +6.39%  '
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("E27").Value = '  +7.16%  '
$ws.Range("E28").Value = '  +12.49%  '
$ws.Range("D29").Value = '0.0₃0738'
$ws.Range("E29").Value = '  +5.41%  '
$ws.Range("E30").Value = '  +4.51%  '
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = '167.07'
$r.ClearFormats()
$ws.Range("E31").Value = '  -2.50%  '
$ws.Range("E32").Value = '  +1.69%  '
$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = '18.38'
$r.ClearFormats()
$ws.Range("E33").Value = '  +2.60%  '
$ws.Range("E34").Value = '  -0.02%  '
$ws.Range("E35").Value = '  -0.24%  '
$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = '1.25'
$r.ClearFormats()
$ws.Range("E36").Value = '  +2.61%  '
$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = '0.917'
$r.ClearFormats()
$ws.Range("E37").Value = '  +1.70%  '
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = '4.00'
$r.ClearFormats()
$ws.Range("E38").Value = '  +4.02%  '
$ws.Range("E39").Value = '  +7.56%  '
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = '37.92'
$r.ClearFormats()
$ws.Range("E40").Value = '  +3.44%  '
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = '0.382'
$r.ClearFormats()
$ws.Range("E41").Value = '  +2.21%  '
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = '140.02'
$r.ClearFormats()
$ws.Range("E42").Value = '  +3.70%  '
$ws.Range("E43").Value = '  +4.38%  '
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = '5.17'
$r.ClearFormats()
$ws.Range("E44").Value = '  +6.53%  '
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = '275.84'
$r.ClearFormats()
$ws.Range("E45").Value = '  +7.66%  '
$ws.Range("E46").Value = '  +2.24%  '
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = '0.0506'
$r.ClearFormats()
$ws.Range("E47").Value = '  +0.59%  '
$ws.Range("E48").Value = '  +2.84%  '
$ws.Range("E49").Value = '  +3.16%  '
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = '0.380'
$r.ClearFormats()
$ws.Range("E50").Value = '  +2.06%  '
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = '17.73'
$r.ClearFormats()
$ws.Range("E51").Value = '  +9.08%  '
